# Demo thong Ke voi HashMap
# Append a new book record ("KD04") as row 5 of the SachDB sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "KD04"
$ws.Range("B5").Value = "NXB_Đại_Học_Quốc_Gia"
$ws.Range("C5").Value = "P.Việt"
$ws.Range("D5").Value = "Thiếu_Nhi"
$ws.Range("E5").Value = "Trường Làng"
$ws.Range("F5").Value = 2020
$ws.Range("G5").Value = 125
$ws.Range("H5").Value = 25000
$ws.Range("I5").Value = "null"

$ws.Columns.Item(2).AutoFit()
